# Adds the strings/rows needed for the new "automatic update check" feature
# to the "Idiomas" (languages) lookup sheet, and rewords the existing
# "problem getting video" message pair.
#
# Final layout of rows 19-22 on the "Idiomas" sheet (A = Spanish, B = English):
#   19  Habido un problema...                    | There was a problem while trying...
#   20  Tienes la ultima versión disponible       | You've the latest version available
#   21  Hay una nueva version disponible...       | There's a new version available...
#   22  Actualizar mas tarde                      | Update later
#
# New shared strings must land in this exact order so the sharedStrings.xml
# unique-string indices line up with the target workbook (A21, B20, B21,
# A20, A22, B22 - matching how the original author apparently typed them).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Idiomas")

# --- Reword row 19 (existing "problem getting video" message pair) -------
$ws.Range("A19").Value = "Habido un problema mientras se intentaba obtener el video.`nPor favor inténtelo de nuevo"
$ws.Range("B19").Value = "There was a problem while trying to get the video.`nPlease try again"

# --- New strings, written in the order that yields the matching shared- --
# --- string table: A21, B20, B21, A20, A22, B22 ---------------------------
$ws.Range("A21").Value = "Hay una nueva version disponible.`n¿Quieres actualizar a la última versión?"
$ws.Range("B20").Value = "You've the latest version available"
$ws.Range("B21").Value = "There's a new version available. `nDo you wanna update to the latest version?"
$ws.Range("A20").Value = "Tienes la ultima versión disponible"
$ws.Range("A22").Value = "Actualizar mas tarde"
$ws.Range("B22").Value = "Update later"

# --- Formatting: rows 20 & 22 use the single-line centered style (like ---
# --- rows 14/16/18); rows 19 & 21 use the centered + wrapped style used --
# --- for the two-line messages (like row 19's original formatting). -----
$ws.Range("A20:B20").HorizontalAlignment = -4108
$ws.Range("A22:B22").HorizontalAlignment = -4108

$ws.Range("A21:B21").HorizontalAlignment = -4108
$ws.Range("A21:B21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 30

# --- Sheet view: scroll position + current selection ----------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B22").Select()
